$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style/formatting (bold, border, centered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Row 2 data
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

# Row 3 data
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7
